$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capitalize first/last names (row 2 = Melvin Leble, row 3 = Pierre Vanobbergen)
$ws.Range("A2").Value = "Melvin"
$ws.Range("B2").Value = "Leble"
$ws.Range("A3").Value = "Pierre"
$ws.Range("B3").Value = "Vanobbergen"

# Update the active selection to I9
$ws.Range("I9").Select() | Out-Null
